$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two time-slot labels in column C (rows 2 and 3)
$ws.Range("C2").Value = "9:30-9:35"
$ws.Range("C3").Value = "9:35-9:40"

# Grow the current selection from C11 to C10:C11
$ws.Range("C11").Select()
$ws.Range("C10:C11").Select()
